# Append: 2025-10-11 06:31 JST
# Update the "取得日時" (retrieved at) timestamp column on the "ランサーズ" sheet
# for all existing data rows (2-7) to the new timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-11 06:31:00"

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
